$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row (Date, Volumen, Precio minimo, Precio maximo, Precio promedio ponderado, Precio $/Kg)
$data = @{
    2  = @{ D = 44362; J = 120; K = 8000; L = 9000; M = 8500; P = 142 }
    3  = @{ D = 44494; J = 120; K = 5000; L = 6000; M = 5500; P = 92 }
    4  = @{ D = 44400; J = 120; K = 9000; L = 10000; M = 9500; P = 158 }
    5  = @{ D = 44421; J = 100; K = 8000; L = 9000; M = 8500; P = 142 }
    7  = @{ D = 44603; J = 140; K = 5500; L = 6000; M = 5750; P = 96 }
    8  = @{ D = 44627; J = 120; K = 4000; L = 4500; M = 4250; P = 71 }
    9  = @{ D = 44648; J = 120; K = 6500; L = 7000; M = 6750; P = 112 }
    10 = @{ D = 44281; J = 120; K = 5500; L = 6000; M = 5750; P = 96 }
    11 = @{ D = 44382; J = 160; K = 7000; L = 8000; M = 7438; P = 124 }
    12 = @{ D = 44589; J = 110; K = 5000; L = 6000; M = 5500; P = 92 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("J$row").Value = $vals.J
    $ws.Range("K$row").Value = $vals.K
    $ws.Range("L$row").Value = $vals.L
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("P$row").Value = $vals.P
}
